$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 141
$ws1.Range("F5").Value = 2987
$ws1.Range("F7").Value = 402

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 141
$ws4.Range("F5").Value = 2987
$ws4.Range("F9").Value = 402
